$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update config paths so they match the new location of the run file
# (files moved under ./src/asset/... instead of ./asset/...)
$ws.Range("I2").Value = "./src/asset/image/example.image.jpg"
$ws.Range("I3").Value = "./src/asset/image/example.image.jpg"
$ws.Range("J2").Value = "./src/asset/font/Sportage-DemoItalic.otf"
$ws.Range("J3").Value = "./src/asset/font/Sportage-DemoItalic.otf"

# Widen column J (10) to fit the longer path text
$ws.Columns.Item(10).ColumnWidth = 35

# Update the selected / active cell shown when the sheet is opened
[void]$ws.Range("I3").Select()
